# Add a new "Change Track" entry row (row 14) below the existing last row (13),
# reusing the formatting of row 13 and appending the new comment required by
# the traceability matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row onto the new row first,
# so the new cells pick up the same styles (date format, centered text, etc.)
# as the rest of the table.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Range("A14").Value = 41755
$ws.Range("B14").Value = "12"
$ws.Range("C14").Value = "JEB"
$ws.Range("D14").Value = "Initial responsibilities asigned in requirements document"
$ws.Range("E14").Value = "Done"

# Move the active selection to the next empty row, matching Excel's usual
# behaviour of advancing selection after data entry.
[void]$ws.Range("A15").Select()
